$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new settings row (row 6) for the "cancel online booking after minute" setting.
# Column order: A = id (uuid), B = name (key), C = desc_name (description), D = value
# Cells are populated in the same order the shared strings were appended in the
# original authoring session: id, name, value, desc_name.
$ws.Range("A6").Value = "771a7a56-6da0-4840-998a-f6e131310dda"
$ws.Range("B6").Value = "cancelOnlineBookingAfterMinute"
$ws.Range("D6").Value = "30"
$ws.Range("C6").Value = "Số phút hủy lịch đặt online"

# Match the selection/active cell left by the author after entering the new row.
$ws.Range("C6").Select()
